$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.35
$ws.Range("C9").Value = -10.614
$ws.Range("A11").Value = -21.032
$ws.Range("B11").Value = 6.504
$ws.Range("A12").Value = -21.465
$ws.Range("C13").Value = -13.112
$ws.Range("C14").Value = -12.25
$ws.Range("A15").Value = -21.767
$ws.Range("C19").Value = -12.303
$ws.Range("C21").Value = -12.213
$ws.Range("C22").Value = -12.34
$ws.Range("B23").Value = 7.392
$ws.Range("C24").Value = -12.681
$ws.Range("C26").Value = -11.845
$ws.Range("A27").Value = -21.506
$ws.Range("A28").Value = -20.648
$ws.Range("B28").Value = 6.341
$ws.Range("A31").Value = -21.858
$ws.Range("A32").Value = -20.841
$ws.Range("B32").Value = 6.572
$ws.Range("B34").Value = 6.901999999999999
$ws.Range("A36").Value = -21.215
$ws.Range("B36").Value = 6.842999999999999
$ws.Range("B37").Value = 6.787999999999999
$ws.Range("A38").Value = -20.513
$ws.Range("C38").Value = -11.941
$ws.Range("C41").Value = -11.842
$ws.Range("B42").Value = 7.761
$ws.Range("A46").Value = -21.547
$ws.Range("B49").Value = 6.748
$ws.Range("C52").Value = -11.552
$ws.Range("A54").Value = -20.908
$ws.Range("B54").Value = 6.092000000000001
$ws.Range("A55").Value = -22.047
$ws.Range("A56").Value = -21.935
$ws.Range("C56").Value = -12.863
$ws.Range("A67").Value = -21.418
$ws.Range("A69").Value = -21.391
$ws.Range("C71").Value = -11.32
$ws.Range("A72").Value = -21.017
$ws.Range("C72").Value = -12.49
$ws.Range("A73").Value = -19.974
$ws.Range("B78").Value = 7.329000000000001
$ws.Range("C78").Value = -11.001
$ws.Range("B80").Value = 7.157999999999999
$ws.Range("A83").Value = -21.117
$ws.Range("C83").Value = -13.183
$ws.Range("C85").Value = -12.346
$ws.Range("A86").Value = -21.298
$ws.Range("C86").Value = -13.233
$ws.Range("C90").Value = -10.614
$ws.Range("A91").Value = -20.813
$ws.Range("A93").Value = -21.327
$ws.Range("C96").Value = -10.311
$ws.Range("B97").Value = 5.449000000000001
$ws.Range("A99").Value = -20.844
$ws.Range("B99").Value = 6.255
$ws.Range("B100").Value = 5.616999999999999
$ws.Range("B101").Value = 5.749000000000001
$ws.Range("C103").Value = -12.869
$ws.Range("A104").Value = -21.164
$ws.Range("A105").Value = -20.451
